$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet (the defined name _xlnm._FilterDatabase updates
# automatically along with the rename since it references this sheet).
$ws.Name = "deneme"

# Update the header row (row 1) text values - these were previously
# formatted like "group[fieldName]" and are now just "fieldName".
$ws.Range("A1").Value = "number"
$ws.Range("B1").Value = "companyName"
$ws.Range("C1").Value = "zipCode"
$ws.Range("D1").Value = "country"
$ws.Range("E1").Value = "city"
$ws.Range("F1").Value = "district"
$ws.Range("G1").Value = "adressDetailText"
$ws.Range("H1").Value = "phone"
$ws.Range("I1").Value = "email"
$ws.Range("J1").Value = "centerType"
$ws.Range("K1").Value = "photo"
$ws.Range("L1").Value = "bio"

# Remove the frozen header pane / split on the sheet view.
$excel.ActiveWindow.FreezePanes = $false

# Return the selection to the natural top-left cell so no stray selection
# state (e.g. the old frozen-pane "B3") lingers in the saved view.
$null = $ws.Range("A1").Select()
